$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-23 07:23:53"
$wsZhCn.Range("G4").Value = "2016-02-23 07:24:42"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-23 07:24:06"
$wsDeDe.Range("G4").Value = "2016-02-23 07:25:05"
